$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix mis-placed "turm" (class) entries so each class appears only once,
# moving it to its correct day/time slot.
$ws.Range("D3").Value = "-"
$ws.Range("E4").Value = "-"
$ws.Range("E8").Value = "EAP"
$ws.Range("C9").Value = "-"
$ws.Range("D10").Value = "EAP"
$ws.Range("E10").Value = "-"
